$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

$nc      = @(20330051920347, 19330051920227, 19330051920239)
$paterno = @("LOPEZ", "CHORA", "RICO")
$materno = @("DE JESUS", "LOPEZ", "ORTIZ")
$nombres = @("EVELYN", "GABRIEL ALEJANDRO", "NADYA GUADALUPE")
$nombreLargo = @(
    "APLICA LA METODOLOGÍA DE DESARROLLO RÁPIDO DE APLICACIONES CON PROGRAMACIÓN ORIENTADA A EVENTOS",
    "CONSTRUYE BASES DE DATOS PARA APLICACIONES WEB",
    "CONSTRUYE BASES DE DATOS PARA APLICACIONES WEB"
)
$grupo   = @("3APM", "5APM", "5APM")
$reprobadas = @(6, 6, 6)

for ($i = 0; $i -lt 3; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $nc[$i]
}
for ($i = 0; $i -lt 3; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $paterno[$i]
}
for ($i = 0; $i -lt 3; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $materno[$i]
}
for ($i = 0; $i -lt 3; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $nombres[$i]
}
for ($i = 0; $i -lt 3; $i++) {
    $ws.Cells.Item($i + 2, 5).Value = $nombreLargo[$i]
}
for ($i = 0; $i -lt 3; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $grupo[$i]
}
for ($i = 0; $i -lt 3; $i++) {
    $ws.Cells.Item($i + 2, 7).Value = $reprobadas[$i]
}
